# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45188 (2023-09-19) to 45189 (2023-09-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
